$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift each date in column A (rows 1-7) forward by 365 days,
# turning the year 2014 into 2015 while keeping month/day/time intact.
for ($r = 1; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 365
}

# Move the sheet's active selection to A9 (was C10).
$ws.Activate()
$ws.Range("A9").Select()
